$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.544.57'
$ws.Range('E2').Value = '  -0.13%  '

$ws.Range('D3').Value = '1.919.23'
$ws.Range('E3').Value = '  -0.27%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4868'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2893'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.77%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06706'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '111.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.26%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.07%  '

$ws.Range('D12').Value = '1.918.46'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07587'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.79%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.279'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6684'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.74%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '293.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.31%  '

$ws.Range('D17').Value = '30.546.87'
$ws.Range('E17').Value = '  -0.22%  '

$ws.Range('E18').Value = '  +0.00%  '

$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007568'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.547'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.79%  '

$ws.Range('D22').Value = '2.169.91'
$ws.Range('E22').Value = '  +0.27%  '

$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.441'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.471'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.70%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.68%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.088'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.79%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1072'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.447'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.26%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.140'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.053'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05017'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7406'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.28%  '

$ws.Range('E35').Value = '  -2.34%  '

$ws.Range('E36').Value = '  +0.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.715'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02028'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.680'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.015'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4430'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.03%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8663'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.65%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.836'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.222'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.69%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.32%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.172'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.59%  '

$ws.Range('E50').Value = '  -0.19%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2533'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.76%  '
